$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistics")

# Remove rows 4-9 (data truncated - emergency break)
$ws.Range("A4:C9").EntireRow.Delete() | Out-Null

# Update row 2 with new timestamp/values
$ws.Range("A2").Value = "2024-08-04 19:40:32"
$ws.Range("B2").Value = 100.1450667798633
$ws.Range("C2").Value = 8

# Update row 3 with new timestamp/values
$ws.Range("A3").Value = "2024-08-04 19:40:34"
$ws.Range("B3").Value = 92.26291639793077
$ws.Range("C3").Value = 16
